$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(30.2640309149665, 11.2035715717842, 4.12540889474795, 55.91849187040492, 3.782953238057835, 10.09795082292669, 12.61300655488907, 22.24781626515347),
    @(29.95707174394685, 10.8124671800033, 4.087986001744706, 55.78947207689728, 3.787673456819107, 10.11795432695479, 12.62176726484598, 22.30910461199847),
    @(29.77591092883957, 10.56899890970743, 4.06441794883955, 55.72390133972421, 3.79071986175247, 10.13093199698501, 12.62931566414995, 22.34875440247665),
    @(29.70400163472665, 10.46912414969471, 4.054666416750521, 55.70061762360822, 3.791998710814223, 10.13639589827031, 12.63293711580206, 22.36541978689812),
    @(29.69217896878287, 10.45250564914057, 4.053038313271044, 55.69695903875726, 3.792213327015807, 10.13731378371964, 12.63357139216176, 22.36821773599539),
    @(29.77493328273244, 10.5676543891708, 4.064287031520503, 55.72357340737727, 3.79073695706168, 10.13100497419891, 12.62936229611686, 22.34897710141293),
    @(30.15672047629354, 11.06951182679648, 4.112628616968339, 55.87117722306186, 3.784550110253329, 10.10470400746658, 12.61557701031959, 22.26852966949716),
    @(30.95989543127078, 12.02016191023425, 4.202699142766219, 56.26855392673846, 3.7735863073466, 10.05862216560988, 12.60575537495846, 22.12677830435717),
    @(31.57825317716452, 12.68973413417336, 4.265937968654899, 56.62560622401765, 3.766233675277597, 10.02808224941008, 12.60902442413133, 22.03237618115013),
    @(31.86459922398702, 12.98657244791185, 4.294058443432807, 56.80195585072071, 3.763039180848184, 10.01490196181499, 12.61278367870404, 21.99154314583947),
    @(31.97366378719698, 13.09775605365212, 4.304612877304749, 56.87071312367525, 3.76185095066794, 10.01001284864995, 12.61453322287781, 21.97638426006222),
    @(31.95014820750337, 13.07386663064691, 4.302343997247664, 56.85581749492123, 3.762105905457339, 10.01106127880622, 12.61414194197358, 21.97963548750652),
    @(31.87355991874625, 12.99574466469976, 4.294928656839318, 56.80757307481549, 3.762940995298945, 10.01449769050741, 12.6129210861209, 21.99028992764999),
    @(31.82672680370306, 12.94773050574278, 4.290374249168627, 56.77827871003716, 3.763455301607391, 10.01661585804856, 12.61221570626749, 21.99685563534171),
    @(31.55963333876215, 12.67016953351071, 4.26408711049726, 56.61435951096978, 3.766445453604689, 10.0289579081157, 12.60882441032876, 22.03508720545705),
    @(31.3970066437869, 12.49782828402387, 4.247794240460212, 56.51735019311024, 3.768318191394898, 10.03671149690365, 12.607325441047, 22.05908179393592),
    @(31.30394931042601, 12.39797795257686, 4.238362161516668, 56.46286536074106, 3.769409492342488, 10.0412382457755, 12.60667719356176, 22.07308152949161),
    @(31.27252723766559, 12.36404945663726, 4.235158204905942, 56.44464383628942, 3.769781423260415, 10.042782462471, 12.60649446873196, 22.07785571956829),
    @(31.41426941359396, 12.51625007736241, 4.249534950307023, 56.52754134607626, 3.768117371707105, 10.03587917408455, 12.60746287378418, 22.05650696519659),
    @(31.89603937344405, 13.01872495638328, 4.297109286820222, 56.8216901758168, 3.762695128073833, 10.01348556928281, 12.61327084033462, 21.98715221666441),
    @(32.21454093545746, 13.33995382790774, 4.327652419570683, 57.02544603562289, 3.759276367766498, 9.999444227709768, 12.61896637707949, 21.94359516222415),
    @(32.04424923750194, 13.16919608222499, 4.311401609777483, 56.91565333153036, 3.761089637757819, 10.00688414567004, 12.6157530260958, 21.96668034261709),
    @(31.40646354102897, 12.5079239717102, 4.248748178043945, 56.52292991489235, 3.768208116679189, 10.0362552519409, 12.60740007529173, 22.05767040697844),
    @(30.73728335182179, 11.76750052393386, 4.178844397780956, 56.1495568246779, 3.776428228706632, 10.07050375496059, 12.60656994293931, 22.16341395584864)
)

$cols = @(2, 3, 4, 6, 7, 10, 12, 14)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($rowIndex, $cols[$j]).Value = $rowData[$j]
    }
}
